$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# The workbook currently ends with a single "总计" (grand total) sheet.
# This edit:
#   1. Turns that existing sheet into the new "2022-Q1" per-fund detail
#      sheet (same shape as the other quarterly sheets).
#   2. Adds a brand-new "总计" sheet (a copy of the original, so it keeps
#      the same look/formatting) with a fresh top row for 2022-Q1 and the
#      existing rows shifted down / re-numbered.
# ---------------------------------------------------------------------------

$totalSheet = $wb.Worksheets.Item("总计")

# Make a copy of the grand-total sheet, placed right after it; this copy
# will become the new grand-total sheet once its data is refreshed.
$totalSheet.Copy($null, $totalSheet)
$newTotalSheet = $wb.Worksheets.Item($totalSheet.Index + 1)
$newTotalSheet.Name = "总计_tmp"

# The original sheet becomes the 2022-Q1 per-fund detail sheet.
$totalSheet.Name = "2022-Q1"

# ---------------------------------------------------------------------------
# Step 1: rebuild "2022-Q1" (ex "总计") as a per-fund detail sheet, matching
# the layout used by 2020-Q4 .. 2021-Q4.
# ---------------------------------------------------------------------------
$q1 = $totalSheet

# Drop the old total rows below the header/data row (rows 3:6).
$q1.Rows("3:6").Delete()

# Extend formatting for the new columns E:H from D1's header style.
$q1.Range("D1").Copy()
$q1.Range("E1:H1").PasteSpecial(-4122)

# Headers.
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Data row (A2 already holds 0 with the right style - leave it as is).
$q1.Range("B2").Value = "'513080"
$q1.Range("C2").Value = "华安法国CAC40ETF（QDII）"
$q1.Range("D2").Value = "'0.60"
$q1.Range("E2").Value = "'96.69"
$q1.Range("F2").Value = "'5.35"
$q1.Range("G2").Value = "'0.0321"
$q1.Range("H2").Value = 5

# ---------------------------------------------------------------------------
# Step 2: rebuild the new "总计" sheet with a 2022-Q1 row inserted on top
# and the remaining rows shifted down / re-numbered.
# ---------------------------------------------------------------------------
$total = $newTotalSheet

$total.Rows("2").Insert()
$total.Range("B2:D2").ClearFormats()

# Copy A column formatting (bold/centered/bordered) onto the new A2 cell.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.03

# Re-number the index column for the rows that shifted down.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5

$total.Name = "总计"

# Restore the originally-active sheet (copying a sheet makes the copy the
# active one, which would otherwise move the workbook's selected tab).
$wb.Worksheets.Item(1).Activate()
